$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 32075
$ws.Range("I21").Value = 300
$ws.Range("K21").Value = 300
$ws.Range("M21").Value = 168
$ws.Range("H23").Value = 32075
$ws.Range("I23").Value = 300
$ws.Range("K23").Value = 300
$ws.Range("M23").Value = -66
$ws.Range("H33").Value = 379.45
$ws.Range("I33").Value = 383.6316
$ws.Range("K33").Value = 383.6316
$ws.Range("M33").Value = -154.6316
$ws.Range("H34").Value = 6098.231
$ws.Range("J34").Value = 23683
$ws.Range("L34").Value = 23683
$ws.Range("N34").Value = -24089
$ws.Range("H36").Value = 6098.231
$ws.Range("J36").Value = 23683
$ws.Range("L36").Value = 23683
$ws.Range("N36").Value = -25113
$ws.Range("H76").Value = 5987.25
$ws.Range("I76").Value = 3299.3333
$ws.Range("J76").Value = 7600
$ws.Range("K76").Value = 3299.3333
$ws.Range("L76").Value = 7600
$ws.Range("M76").Value = -2984.3333
$ws.Range("N76").Value = -8230
$ws.Range("H79").Value = 5987.25
$ws.Range("I79").Value = 3299.3333
$ws.Range("J79").Value = 7600
$ws.Range("K79").Value = 3299.3333
$ws.Range("L79").Value = 7600
$ws.Range("M79").Value = -2207.3333
$ws.Range("N79").Value = -9784
$ws.Range("H80").Value = 1653.6
$ws.Range("I80").Value = 2926
$ws.Range("J80").Value = 805.3333
$ws.Range("K80").Value = 8778
$ws.Range("L80").Value = 2415.9999
$ws.Range("M80").Value = -7780
$ws.Range("N80").Value = -4411.9999
$ws.Range("H83").Value = 1653.6
$ws.Range("I83").Value = 2926
$ws.Range("J83").Value = 805.3333
$ws.Range("K83").Value = 26334
$ws.Range("L83").Value = 7247.9997
$ws.Range("M83").Value = -21342
$ws.Range("N83").Value = -17231.9997
$ws.Range("H87").Value = 57687.332
$ws.Range("J87").Value = 57687.332
$ws.Range("L87").Value = 57687.332
$ws.Range("N87").Value = -60183.332
$ws.Range("H90").Value = 57687.332
$ws.Range("J90").Value = 57687.332
$ws.Range("L90").Value = 173061.996
$ws.Range("N90").Value = -185541.996
$ws.Range("H129").Value = 1205.0426
$ws.Range("I129").Value = 493.7
$ws.Range("J129").Value = 1397.2972
$ws.Range("K129").Value = 1481.1
$ws.Range("L129").Value = 4191.8916
$ws.Range("M129").Value = 3518.9
$ws.Range("N129").Value = -14191.8916
$ws.Range("H139").Value = 200000
$ws.Range("J139").Value = 200000
$ws.Range("L139").Value = 200000
$ws.Range("N139").Value = -210280

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8429.885
$ws.Range("I32").Value = 8807.308999999999
$ws.Range("J32").Value = 7079.1055
$ws.Range("K32").Value = 8807.308999999999
$ws.Range("L32").Value = 7079.1055
$ws.Range("M32").Value = -8520.308999999999
$ws.Range("N32").Value = -7653.1055
$ws.Range("H132").Value = 5320960
$ws.Range("I132").Value = 5954002
$ws.Range("J132").Value = 3405.6
$ws.Range("K132").Value = 17862006
$ws.Range("L132").Value = 10216.8
$ws.Range("M132").Value = -17859476
$ws.Range("N132").Value = -15276.8

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 614.2857
$ws.Range("I22").Value = 380
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -30
$ws.Range("N22").Value = -1900

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 534251.5600000001
$ws.Range("I2").Value = 96.46154
$ws.Range("J2").Value = 1068406.6
$ws.Range("K2").Value = 578.76924
$ws.Range("L2").Value = 6410439.600000001
$ws.Range("M2").Value = -465.76924
$ws.Range("N2").Value = -6410665.600000001
$ws.Range("H38").Value = 236.3077
$ws.Range("I38").Value = 831
$ws.Range("J38").Value = 128.18182
$ws.Range("K38").Value = 2493
$ws.Range("L38").Value = 384.5454599999999
$ws.Range("M38").Value = -2146
$ws.Range("N38").Value = -1078.54546
$ws.Range("H107").Value = 1151.1154
$ws.Range("I107").Value = 139.8
$ws.Range("J107").Value = 1391.9048
$ws.Range("K107").Value = 419.4
$ws.Range("L107").Value = 4175.7144
$ws.Range("M107").Value = 1500.6
$ws.Range("N107").Value = -8015.7144
$ws.Range("H132").Value = 1795.4445
$ws.Range("I132").Value = 751
$ws.Range("J132").Value = 2631
$ws.Range("K132").Value = 6759
$ws.Range("L132").Value = 23679
$ws.Range("M132").Value = -4229
$ws.Range("N132").Value = -28739

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11907972
$ws.Range("I80").Value = 20835816
$ws.Range("J80").Value = 2384939.8
$ws.Range("K80").Value = 20835816
$ws.Range("L80").Value = 2384939.8
$ws.Range("M80").Value = -20834818
$ws.Range("N80").Value = -2386935.8
$ws.Range("H83").Value = 11907972
$ws.Range("I83").Value = 20835816
$ws.Range("J83").Value = 2384939.8
$ws.Range("K83").Value = 104179080
$ws.Range("L83").Value = 11924699
$ws.Range("M83").Value = -104174088
$ws.Range("N83").Value = -11934683

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4920.7607
$ws.Range("I7").Value = 5162.9546
$ws.Range("J7").Value = 4698.75
$ws.Range("K7").Value = 5162.9546
$ws.Range("L7").Value = 4698.75
$ws.Range("M7").Value = -5050.9546
$ws.Range("N7").Value = -4922.75
$ws.Range("H126").Value = 4920.7607
$ws.Range("I126").Value = 5162.9546
$ws.Range("J126").Value = 4698.75
$ws.Range("K126").Value = 15488.8638
$ws.Range("L126").Value = 14096.25
$ws.Range("M126").Value = -13018.8638
$ws.Range("N126").Value = -19036.25
$ws.Range("H132").Value = 15637005
$ws.Range("I132").Value = 9342.714
$ws.Range("J132").Value = 27791854
$ws.Range("K132").Value = 28028.142
$ws.Range("L132").Value = 83375562
$ws.Range("M132").Value = -25498.142
$ws.Range("N132").Value = -83380622
$ws.Range("H139").Value = 55716.668
$ws.Range("J139").Value = 55716.668
$ws.Range("L139").Value = 55716.668
$ws.Range("N139").Value = -65996.66800000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16257.429
$ws.Range("I62").Value = 6966.3335
$ws.Range("J62").Value = 23225.75
$ws.Range("K62").Value = 6966.3335
$ws.Range("L62").Value = 23225.75
$ws.Range("M62").Value = -6342.3335
$ws.Range("N62").Value = -24473.75
$ws.Range("H65").Value = 16257.429
$ws.Range("I65").Value = 6966.3335
$ws.Range("J65").Value = 23225.75
$ws.Range("K65").Value = 34831.6675
$ws.Range("L65").Value = 116128.75
$ws.Range("M65").Value = -31711.6675
$ws.Range("N65").Value = -122368.75
$ws.Range("H81").Value = 1053.3334
$ws.Range("I81").Value = 433.33334
$ws.Range("J81").Value = 1260
$ws.Range("K81").Value = 866.66668
$ws.Range("L81").Value = 2520
$ws.Range("M81").Value = 194.33332
$ws.Range("N81").Value = -4642
$ws.Range("H84").Value = 1053.3334
$ws.Range("I84").Value = 433.33334
$ws.Range("J84").Value = 1260
$ws.Range("K84").Value = 4333.3334
$ws.Range("L84").Value = 12600
$ws.Range("M84").Value = 970.6665999999996
$ws.Range("N84").Value = -23208
$ws.Range("H136").Value = 991.2195
$ws.Range("I136").Value = 776.25714
$ws.Range("J136").Value = 2245.1667
$ws.Range("K136").Value = 2328.77142
$ws.Range("L136").Value = 6735.500100000001
$ws.Range("M136").Value = 221.22858
$ws.Range("N136").Value = -11835.5001

Write-Output "Applied 190 cell updates"